# Fix duplicate duties of "Night Tests" duty roster, and tidy up the
# stale selection left on the "Demo" sheet.

$wb = $excel.ActiveWorkbook

# --- "Demo" sheet: clean up the saved selection (was a stray E10 cell) ---
$wsDemo = $wb.Worksheets.Item("Demo")
$wsDemo.Activate()
[void]$wsDemo.Range("A1:A4").Select()

# --- "Night Tests" sheet: re-roster the duty list so nobody repeats back ---
$wsNight = $wb.Worksheets.Item("Night Tests")
$wsNight.Activate()

# New rotation for rows 1-4, 6-9, 11-14, 16-19 (in place) plus new rows
# 21-24; rows 5, 10, 15, 20 are cleared out entirely (no more blank-row
# placeholders sitting mid-list with leftover duplicate names).
$wsNight.Range("A1").Value = "Andrii Vanikhin"
$wsNight.Range("A2").Value = "Eugene Zinchenko"
$wsNight.Range("A3").Value = "Olena Mikheyeva"
$wsNight.Range("A4").Value = "Dmytro Latyshko"
$wsNight.Range("A5").Clear()

$wsNight.Range("A6").Value = "Eugene Zinchenko"
$wsNight.Range("A7").Value = "Dmytro Latyshko"
$wsNight.Range("A8").Value = "Olena Mikheyeva"
$wsNight.Range("A9").Value = "Andrii Vanikhin"
$wsNight.Range("A10").Clear()

$wsNight.Range("A11").Value = "Olena Mikheyeva"
$wsNight.Range("A12").Value = "Andrii Vanikhin"
$wsNight.Range("A13").Value = "Eugene Zinchenko"
$wsNight.Range("A14").Value = "Dmytro Latyshko"
$wsNight.Range("A15").Clear()

$wsNight.Range("A16").Value = "Eugene Zinchenko"
$wsNight.Range("A17").Value = "Olena Mikheyeva"
$wsNight.Range("A18").Value = "Andrii Vanikhin"
$wsNight.Range("A19").Value = "Dmytro Latyshko"
$wsNight.Range("A20").Clear()

$wsNight.Range("A21").Value = "Andrii Vanikhin"
$wsNight.Range("A22").Value = "Eugene Zinchenko"
$wsNight.Range("A23").Value = "Dmytro Latyshko"
$wsNight.Range("A24").Value = "Olena Mikheyeva"

[void]$wsNight.Range("A1:A24").Select()
